$wb = $excel.ActiveWorkbook

# Sheets with Ano-prefixed year headers (B1:E1 = 2015/2030/2040/2050)
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet with Intervalo-prefixed headers (B1:E1 = 2015/2015-2030/2031-2040/2041-2050)
$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIntervalo.Range("B1").Value = "Intervalo 2015"
$wsIntervalo.Range("C1").Value = "Intervalo 2015-2030"
$wsIntervalo.Range("D1").Value = "Intervalo 2031-2040"
$wsIntervalo.Range("E1").Value = "Intervalo 2041-2050"

# Sheet with only a single year header (B1 = 2015)
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano 2015"
